# Update the "Förändrad" date column (C) for rows 2-27 from 2023-09-23 (45192)
# to 2023-10-03 (45202), keeping the existing number format/style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newSerial = 45202

for ($row = 2; $row -le 27; $row++) {
    $cell = $ws.Cells.Item($row, 3)   # Column C
    if ($cell.Value2 -eq 45192) {
        $cell.Value = $newSerial
    }
}
